$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 304, shifting existing rows 304:337 down to 305:338
$ws.Rows.Item(304).Insert()

# Populate the newly inserted row 304 with the new weekly record
$ws.Cells.Item(304, 1).Value = 11
$ws.Cells.Item(304, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(304, 3).Value = "Bíobío"
$ws.Cells.Item(304, 4).Value = 45154
$ws.Cells.Item(304, 5).Value = 8
$ws.Cells.Item(304, 6).Value = "Fruta"
$ws.Cells.Item(304, 7).Value = 100101
$ws.Cells.Item(304, 8).Value = "Berries"
$ws.Cells.Item(304, 9).Value = 100101007
$ws.Cells.Item(304, 10).Value = "Kiwi"
$ws.Cells.Item(304, 11).Value = "Hayward"
$ws.Cells.Item(304, 12).Value = "Primera"
$ws.Cells.Item(304, 13).Value = 100
$ws.Cells.Item(304, 14).Value = 15000
$ws.Cells.Item(304, 15).Value = 15000
$ws.Cells.Item(304, 16).Value = 15000
$ws.Cells.Item(304, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(304, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(304, 19).Value = 833
$ws.Cells.Item(304, 20).Value = 18
